# Applies:
#  1. Insert a new slide ("역할" / Roles) as the 2nd slide, duplicating the
#     "Title and body" layout/formatting used by the existing "사업 개요" slide.
#  2. On the slide that ends up last ("확장 방향성" / Expansion direction),
#     merge the final two text runs of the last body paragraph into one run
#     (no visible text change - purely a run-structure normalization that
#     happens when the paragraph text is re-set as a whole).

$p = $ppt.ActivePresentation

# --- 1. Insert the new "역할" slide at position 2 -----------------------
# Duplicate slide 2 ("사업 개요") so the new slide inherits the same
# slide layout (Title and body) and shape/placeholder structure, then move
# the duplicate into position 2 and overwrite its text.
$srcSlide = $p.Slides.Item(2)
$dupRange = $srcSlide.Duplicate()
$newSlide = $dupRange.Item(1)
$newSlide.MoveTo(2)

# Title placeholder
$titleShape = $newSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "역할"

# Body placeholder - one paragraph per teammate, blank paragraph between
$bodyShape = $newSlide.Shapes.Item(2)
$bodyText = "박해온- 코드 개발, 사진 촬영, 초기 아이디어 제공`r`r심민준- 코드 개발, 연기, 아두이노 조립`r`r이동하- 사업 가능성 및 확장 방향성 아이디어 제공`r`r최윤후- 코드 개발, 연기, 사업가능성,확장 방향성 아이디어 제공"
$bodyShape.TextFrame.TextRange.Text = $bodyText

# --- 2. Normalize the run split on the last slide's last paragraph ------
# Find the "확장 방향성" slide (now pushed one position later than before).
$lastSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.Shapes.Item(1).TextFrame.TextRange.Text -eq "확장 방향성") {
        $lastSlide = $candidate
    }
}

if ($lastSlide -ne $null) {
    $expShape = $lastSlide.Shapes.Item(2)
    $expTr = $expShape.TextFrame.TextRange
    $paraCount = $expTr.Paragraphs().Count
    $lastPara = $expTr.Paragraphs($paraCount, 1)

    # Locate "정신과" (the start of the run that should absorb the trailing
    # "수 있음" run) within the paragraph text.
    $marker = "정신과 의사와의 협력을 통해 제품을 홍보하고 신뢰성을 높일 수 있음"
    $paraText = $lastPara.Text
    $relIdx = $paraText.IndexOf("정신과 의사와의")
    if ($relIdx -ge 0) {
        $absStart = $lastPara.Start + $relIdx
        $tailLen = $paraText.Length - $relIdx
        $target = $expTr.Characters($absStart, $tailLen)
        # Overwrite via a distinct placeholder first so the underlying
        # engine always rewrites the run (no-op skip when text is unchanged).
        $target.Text = "placeholder_tmp"
        $tr2 = $expShape.TextFrame.TextRange
        $target2 = $tr2.Characters($absStart, "placeholder_tmp".Length)
        $target2.Text = $marker
    }
}
